# Fix the "Company Name" merge-field placeholder paragraph in the
# observation-form table: correct the "compayNameForOF" typo to
# "companyNameForOF", refresh its run formatting (color/size/lang),
# mark the surrounding text with grammar proof-marks, and append the
# trailing blank-space runs used for printing extra room after the
# field (temperature / RH spacing), per the template update.

$d = $word.ActiveDocument

# Locate the exact paragraph that still needs fixing so we can carry
# over its original paragraph-mark attributes (paraId/rsid/etc.) and
# paragraph properties (<w:pPr>) unchanged.
$full = $d.Content.WordOpenXML
$needle = "compayNameForOF"
$idx = $full.IndexOf($needle)
if ($idx -lt 0) {
    throw "Could not locate 'compayNameForOF' placeholder text in document"
}

$before = $full.Substring(0, $idx)
$pStart = $before.LastIndexOf("<w:p ")
if ($pStart -lt 0) {
    throw "Could not locate enclosing <w:p> for placeholder"
}

$tagEnd = $full.IndexOf(">", $pStart)
$pOpenTag = $full.Substring($pStart, $tagEnd - $pStart + 1)

$afterOpenTag = $full.Substring($tagEnd + 1)
$pPr = ""
if ($afterOpenTag.StartsWith("<w:pPr>")) {
    $pPrEnd = $afterOpenTag.IndexOf("</w:pPr>") + "</w:pPr>".Length
    $pPr = $afterOpenTag.Substring(0, $pPrEnd)
}

# Select the whole "{compayNameForOF}" placeholder run-content (brace
# to brace) so InsertXML replaces exactly that span while leaving the
# table cell / row / paragraph mark itself untouched.
$rng = $d.Content
$found = $rng.Find.Execute("{compayNameForOF}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find placeholder range to replace"
}

# New run content: corrected spelling, refreshed rPr, grammar
# proof-marks, and the extra trailing spacer runs added after the
# closing brace.
$newRuns = '<w:r><w:rPr><w:szCs w:val="32"/></w:rPr><w:t>{</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>companyNameForOF</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:szCs w:val="32"/></w:rPr><w:t>}</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="000000"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t xml:space="preserve">   </w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:color w:val="000000"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t xml:space="preserve">                    </w:t></w:r>'

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
    '<w:body>' + $pOpenTag + $pPr + $newRuns + '</w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

[void]$rng.InsertXML($xmlFrag)

Write-Output "Placeholder paragraph updated successfully"
